$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.419591
$ws.Cells.Item(2, 8).Value = 4.258773
$ws.Cells.Item(2, 9).Value = 0.001848767113890483
$ws.Cells.Item(2, 10).Value = 0.001848767113890483
$ws.Cells.Item(2, 13).Value = 8.142376000000001
$ws.Cells.Item(2, 14).Value = 24.427128
$ws.Cells.Item(2, 15).Value = 0.1741313933276368
$ws.Cells.Item(2, 16).Value = 0.1741313933276368
$ws.Cells.Item(2, 17).Value = 11.558843688216
$ws.Cells.Item(2, 18).Value = 104.029593193944
$ws.Cells.Item(2, 19).Value = 0.0003219283934800636
$ws.Cells.Item(2, 20).Value = 0.0003219283934800636
$ws.Cells.Item(3, 7).Value = 1.419591
$ws.Cells.Item(3, 8).Value = 4.258773
$ws.Cells.Item(3, 9).Value = 0.001848767113890483
$ws.Cells.Item(3, 10).Value = 0.001848767113890483
$ws.Cells.Item(3, 15).Value = 0.5205382400466131
$ws.Cells.Item(3, 16).Value = 0.5205382400466131
$ws.Cells.Item(3, 17).Value = 34.553333752501
$ws.Cells.Item(3, 18).Value = 310.980003772509
$ws.Cells.Item(3, 19).Value = 0.0009623539797206085
$ws.Cells.Item(3, 20).Value = 0.0009623539797206085
$ws.Cells.Item(4, 7).Value = 1.419591
$ws.Cells.Item(4, 8).Value = 4.258773
$ws.Cells.Item(4, 9).Value = 0.001848767113890483
$ws.Cells.Item(4, 10).Value = 0.001848767113890483
$ws.Cells.Item(4, 15).Value = 0.3053303666257501
$ws.Cells.Item(4, 16).Value = 0.3053303666257501
$ws.Cells.Item(4, 17).Value = 20.267832891294
$ws.Cells.Item(4, 18).Value = 182.410496021646
$ws.Cells.Item(4, 19).Value = 0.0005644847406898112
$ws.Cells.Item(4, 20).Value = 0.0005644847406898112
$ws.Cells.Item(5, 9).Value = 0.9578582377148513
$ws.Cells.Item(5, 10).Value = 0.9578582377148513
$ws.Cells.Item(5, 13).Value = 8.142376000000001
$ws.Cells.Item(5, 14).Value = 24.427128
$ws.Cells.Item(5, 15).Value = 0.1741313933276368
$ws.Cells.Item(5, 16).Value = 0.1741313933276368
$ws.Cells.Item(5, 17).Value = 5988.711916189934
$ws.Cells.Item(5, 18).Value = 53898.40724570941
$ws.Cells.Item(5, 19).Value = 0.1667931895436418
$ws.Cells.Item(5, 20).Value = 0.1667931895436418
$ws.Cells.Item(6, 9).Value = 0.9578582377148513
$ws.Cells.Item(6, 10).Value = 0.9578582377148513
$ws.Cells.Item(6, 15).Value = 0.5205382400466131
$ws.Cells.Item(6, 16).Value = 0.5205382400466131
$ws.Cells.Item(6, 19).Value = 0.4986018412742391
$ws.Cells.Item(6, 20).Value = 0.4986018412742391
$ws.Cells.Item(7, 9).Value = 0.9578582377148513
$ws.Cells.Item(7, 10).Value = 0.9578582377148513
$ws.Cells.Item(7, 15).Value = 0.3053303666257501
$ws.Cells.Item(7, 16).Value = 0.3053303666257501
$ws.Cells.Item(7, 19).Value = 0.2924632068969705
$ws.Cells.Item(7, 20).Value = 0.2924632068969705
$ws.Cells.Item(8, 8).Value = 92.81792100000001
$ws.Cells.Item(8, 9).Value = 0.04029299517125823
$ws.Cells.Item(8, 10).Value = 0.04029299517125823
$ws.Cells.Item(8, 13).Value = 8.142376000000001
$ws.Cells.Item(8, 14).Value = 24.427128
$ws.Cells.Item(8, 15).Value = 0.1741313933276368
$ws.Cells.Item(8, 16).Value = 0.1741313933276368
$ws.Cells.Item(8, 17).Value = 251.919470773432
$ws.Cells.Item(8, 18).Value = 2267.275236960889
$ws.Cells.Item(8, 19).Value = 0.007016275390514936
$ws.Cells.Item(8, 20).Value = 0.007016275390514936
$ws.Cells.Item(9, 8).Value = 92.81792100000001
$ws.Cells.Item(9, 9).Value = 0.04029299517125823
$ws.Cells.Item(9, 10).Value = 0.04029299517125823
$ws.Cells.Item(9, 15).Value = 0.5205382400466131
$ws.Cells.Item(9, 16).Value = 0.5205382400466131
$ws.Cells.Item(9, 17).Value = 753.0733858147105
$ws.Cells.Item(9, 19).Value = 0.02097404479265344
$ws.Cells.Item(9, 20).Value = 0.02097404479265344
$ws.Cells.Item(10, 8).Value = 92.81792100000001
$ws.Cells.Item(10, 9).Value = 0.04029299517125823
$ws.Cells.Item(10, 10).Value = 0.04029299517125823
$ws.Cells.Item(10, 15).Value = 0.3053303666257501
$ws.Cells.Item(10, 16).Value = 0.3053303666257501
$ws.Cells.Item(10, 19).Value = 0.01230267498808985
$ws.Cells.Item(10, 20).Value = 0.01230267498808985

$wb.Save()
